$d = $word.ActiveDocument

# The "Requisitos" section ends with the paragraph
# "LOQ4010: Introdução à Engenharia Química (Requisito fraco)".
# Immediately after it, the document used to contain three more
# paragraphs (leftover page-footer boilerplate from the site build):
#   1) an empty paragraph
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and
#      Github pages. Original theme under Creative Commons Attribution"
# Those three paragraphs are removed, leaving the LOQ4010 paragraph
# followed directly by the pre-existing blank paragraph / page break.

$anchorText = "LOQ4010: Introdução à Engenharia Química (Requisito fraco)"

$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd("`r", "`n") -eq $anchorText) {
        $anchor = $p
        break
    }
}

if ($anchor -ne $null) {
    $toRemove = New-Object System.Collections.ArrayList

    $cur = $anchor.Next()
    for ($k = 0; $k -lt 3; $k++) {
        $toRemove.Add($cur) | Out-Null
        $cur = $cur.Next()
    }

    # Delete furthest paragraph first so the earlier ones keep valid
    # ranges.
    for ($k = $toRemove.Count - 1; $k -ge 0; $k--) {
        $toRemove[$k].Range.Delete()
    }
}
